# Qatar Stars League update — swap the content (columns B..AC) of the
# given row pairs while leaving column A (the sequential index) untouched.
# This corresponds to the upstream "re-sorting" of rows seen in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(15, 16),
    @(24, 25),
    @(29, 30),
    @(50, 51),
    @(56, 57),
    @(60, 61),
    @(68, 69),
    @(78, 79),
    @(81, 82),
    @(87, 88),
    @(98, 99),
    @(104, 105),
    @(106, 107),
    @(108, 109),
    @(110, 111),
    @(114, 115),
    @(118, 119),
    @(120, 121)
)

# Columns B (2) through AC (29) inclusive.
$firstCol = 2
$lastCol = 29

foreach ($pair in $pairs) {
    $rowA = $pair[0]
    $rowB = $pair[1]

    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cellA = $ws.Cells.Item($rowA, $col)
        $cellB = $ws.Cells.Item($rowB, $col)

        $valA = $cellA.Value2
        $valB = $cellB.Value2

        $cellA.Value2 = $valB
        $cellB.Value2 = $valA
    }
}
